# Refactored insert_data_to_db: batch execution replaces the old
# duplicate-check + row-by-row insert flow. Reflect the new data shape
# in the "Demo_Excel" report sheet:
#   - Pages Printed is now a real accumulated count (no longer a fixed
#     per-row value of 1 / 10)
#   - Date is now written as a genuine Excel date/time value (was text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (User1 / land AC.xlsx / Host1) ---
$ws.Range("A2").Value = "User1"
$ws.Range("B2").Value = "land AC.xlsx"
$ws.Range("C2").Value = "Host1"
$ws.Range("D2").Value = 1000
$ws.Range("E2").NumberFormat = "m/d/yy h:mm"
$ws.Range("E2").Value = 46026.023622685185

# --- Row 3 (User2 / Shift Reporting format.xlsx / Host2) ---
$ws.Range("A3").Value = "User2"
$ws.Range("B3").Value = "Shift Reporting format.xlsx"
$ws.Range("C3").Value = "Host2"
$ws.Range("D3").Value = 50000
$ws.Range("E3").NumberFormat = "m/d/yy h:mm"
$ws.Range("E3").Value = 46026.026006944441

# Row 2 no longer needs the taller wrap height now that the date column
# holds a short formatted date/time instead of wrapped text.
$ws.Rows.Item(2).RowHeight = 16

# Move the active selection as left by the author's last save.
$ws.Range("E5").Select()
